$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells to keep their exact string values (these
# columns hold inline-string data in the original workbook, including
# numeric-looking text like "215.15" that Excel would otherwise coerce
# to a genuine number on assignment).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "89.226.31"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -3.75%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.136.28"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -4.30%  "
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.15"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -1.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "635.92"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +0.89%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.395"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -1.90%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.772"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +8.13%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.135.26"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -4.24%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.561"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -4.53%  "
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -0.57%  "
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -6.22%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.33"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -0.31%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "89.126.43"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -3.55%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.708.48"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -4.43%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "32.35"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -5.61%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.133.65"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -4.03%  "
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "SuiNetwork"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.39"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +2.91%  "
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "PEPE"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000228"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +4.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.27"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -5.05%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "427.22"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -4.86%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.37"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -5.33%  "
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -6.38%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.50"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +2.75%  "
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +6.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.58"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -4.45%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "3.294.80"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -4.57%  "
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = "Cronos"
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.157"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -10.35%  "
$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = "Binance-PegBSC-USD"
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.976"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -2.32%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +4.32%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "8.22"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -6.11%  "
$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = "Kaspa"
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +16.61%  "
$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = "Bittensor"
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "507.11"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -8.98%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "7.08"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -0.90%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.29"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +0.43%  "
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -4.79%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "22.05"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -2.57%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "22.24"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -0.96%  "
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.87"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -5.82%  "
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -7.12%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "145.79"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -2.95%  "
$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.131"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +2.79%  "
$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = "OKB"
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "43.77"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -3.62%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "164.25"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -8.38%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0648"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +7.65%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.725"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -1.02%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "24.24"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -2.98%  "
